$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 9.5

$ws.Range("C5").Select()
